$d = $word.ActiveDocument

# Locate the paragraph that currently uses the "ResourceHeadnote" style.
# It holds a single, empty run ("<w:r><w:t xml:space=\"preserve\"/></w:r>").
# The target state turns it into two separate, completely run-less
# paragraphs that both use the "CaseText" style.
#
# Just reassigning .Style (or clearing .Range.Text) leaves an (now
# empty) run behind, so instead we replace the paragraph's range with
# raw OOXML describing exactly the two run-less paragraphs we want.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Resource Headnote") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    $target = $d.Paragraphs.Item(3)
}

$range = $target.Range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p><w:pPr><w:pStyle w:val="CaseText"/></w:pPr></w:p>' +
       '<w:p><w:pPr><w:pStyle w:val="CaseText"/></w:pPr></w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData>' +
       '</pkg:part>' +
       '</pkg:package>'

$range.InsertXML($xml)
